# Update odds-data rows for "Ecuador LigaPro Serie A".
#
# The underlying source data got re-synced and a handful of match rows
# ended up shuffled relative to their original row position (the id in
# column B, the two team names, score and every odds column move together
# as one record, while the running row index in column A and the
# Div/Date columns C:D stay where they are). This script reads the
# existing per-row data for the affected rows, then writes it back out
# in the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that belong to a single match record, excluding A (row
# index), C (Div) and D (Date) which never change here.
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$row").Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $val = $data[$c]
        if ($val -ne $null) {
            $ws.Range("$c$row").Value = $val
        }
    }
}

# --- snapshot current contents of every row that is going to move -------
$rows1 = @(142, 143, 144, 145)
$snap1 = @{}
foreach ($r in $rows1) { $snap1[$r] = Get-RowData $r }

$rows2 = @(254, 255, 256)
$snap2 = @{}
foreach ($r in $rows2) { $snap2[$r] = Get-RowData $r }

# --- apply the new row order ---------------------------------------------
# Group 1 (rows 142-145): row 142 receives the old row-145 record,
# row 143 receives the old row-144 record, row 144 receives the old
# row-142 record and row 145 receives the old row-143 record.
$map1 = @{142 = 145; 143 = 144; 144 = 142; 145 = 143}
foreach ($newRow in $rows1) {
    Set-RowData $newRow $snap1[$map1[$newRow]]
}

# Group 2 (rows 254-256): row 254 receives the old row-256 record,
# row 255 receives the old row-254 record and row 256 receives the old
# row-255 record.
$map2 = @{254 = 256; 255 = 254; 256 = 255}
foreach ($newRow in $rows2) {
    Set-RowData $newRow $snap2[$map2[$newRow]]
}
